$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (Tabelle1 -> ps1)
$ws.Name = "ps1"

# Translate header row
$ws.Range("B2").Value = "Unit"
$ws.Range("C2").Value = "Pump storage 1"

# Translate parameter names (column A) and units (column B)
$ws.Range("A3").Value = "Turbine power (net) max"
$ws.Range("B3").Value = "MW"

$ws.Range("A4").Value = "Turbine power (net)min"
$ws.Range("B4").Value = "MW"

$ws.Range("A5").Value = "Pump power (gross) max"
$ws.Range("B5").Value = "MW"

$ws.Range("A6").Value = "Pumpleistung (gross) min"
$ws.Range("B6").Value = "MW"

$ws.Range("A7").Value = "Turbine efficiency"

$ws.Range("A8").Value = "Pump efficiency"

$ws.Range("A9").Value = "Net energy capacity"
